$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("invalid")

# The "credit_purpose" (H) and "credit_purpose_ff" (I) columns were part of an
# abandoned merge and need to go away entirely; delete both columns outright
# (their custom widths go with them, and whatever used to sit to the right of
# them - already at the sheet's default width - slides in to take their place).
$ws.Columns("H:I").Delete()

# Column H (now empty/default) becomes the new "action_taken" column.
$ws.Range("H1").Value = "action_taken"

# Fill in the action_taken values row by row (plain numbers, no special type).
$ws.Range("H2").Value = 1
$ws.Range("H3").Value = 2
$ws.Range("H4").Value = 3
$ws.Range("H5").Value = 4
$ws.Range("H6").Value = 5
$ws.Range("H7").Value = 6
$ws.Range("H8").Value = 1
$ws.Range("H9").Value = 2
$ws.Range("H10").Value = 3
# Row 11 is intentionally left blank in the action_taken column.

# Fix up the ct_credit_product values that were wrong after the merge.
$ws.Range("D10").Value = 977
$ws.Range("D11").Value = 988

# The header row got taller (user widened it while editing).
$ws.Rows.Item(1).RowHeight = 34

# Leave the selection/view where the editor ended up.
$ws.Activate() | Out-Null
$ws.Range("H11").Select() | Out-Null
